$p = $ppt.ActivePresentation

# --- 1) Table style swap on the three "balance sheet" tables (slides 14-16) ---
foreach ($idx in 14, 15, 16) {
    $s = $p.Slides.Item($idx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{248D38CB-E44F-479A-9FAB-0652DC57B296}", $true)
        }
    }
}

# --- 2) Theme colour swap: the deck's active theme goes from "Integral"
#        (Red Violet) to the plain default "Office Theme" palette. ---
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Colors(1).RGB  = 0x00000000  # dk1     000000
$tcs.Colors(2).RGB  = 0x00FFFFFF  # lt1     FFFFFF
$tcs.Colors(3).RGB  = 0x006A5444  # dk2     44546A
$tcs.Colors(4).RGB  = 0x00E6E6E7  # lt2     E7E6E6
$tcs.Colors(5).RGB  = 0x00D59B5B  # accent1 5B9BD5
$tcs.Colors(6).RGB  = 0x00317DED  # accent2 ED7D31
$tcs.Colors(7).RGB  = 0x00A5A5A5  # accent3 A5A5A5
$tcs.Colors(8).RGB  = 0x0000C0FF  # accent4 FFC000
$tcs.Colors(9).RGB  = 0x00C47244  # accent5 4472C4
$tcs.Colors(10).RGB = 0x0047AD70  # accent6 70AD47
$tcs.Colors(11).RGB = 0x00C16305  # hlink   0563C1
$tcs.Colors(12).RGB = 0x00724F95  # folHlink 954F72
